$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule name for the 4th rule row (R40 -> "1") in the "Rules String Hello"
# decision table. Use an apostrophe (text) prefix so the numeric-looking
# value "1" is stored as text, matching the existing column's string type.
$ws.Range("B11").Value = "'1"
